# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds dotted numeric-looking strings (e.g. "26.259.30")
# that must stay literal text -- force text format before assigning so the
# COM layer does not coerce them into doubles (which would mangle the digits
# and drop trailing zeros).
$priceCells = @(
    "D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D20", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D44", "D46", "D48", "D50", "D51"
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.259.30"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "1.689.14"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "219.02"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").Value = "0.5246"
$ws.Range("E6").Value = "  +3.45%  "
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "0.2695"
$ws.Range("E8").Value = "  +2.25%  "
$ws.Range("B9").Value = "Solana"
$ws.Range("C9").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D9").Value = "22.23"
$ws.Range("E9").Value = "  +2.93%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.06454"
$ws.Range("E10").Value = "  +2.14%  "
$ws.Range("D11").Value = "0.07466"
$ws.Range("E11").Value = "  +0.99%  "
$ws.Range("D12").Value = "1.688.70"
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("D13").Value = "4.554"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").Value = "0.5875"
$ws.Range("E14").Value = "  +2.07%  "
$ws.Range("D15").Value = "0.000008589"
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D16").Value = "64.93"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").Value = "26.319.18"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "4.989"
$ws.Range("E18").Value = "  +0.60%  "
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("D20").Value = "10.88"
$ws.Range("E20").Value = "  +0.78%  "
$ws.Range("D21").Value = "191.28"
$ws.Range("E21").Value = "  +2.30%  "
$ws.Range("D22").Value = "6.253"
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "145.33"
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "7.678"
$ws.Range("E25").Value = "  +0.77%  "
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").Value = "0.1243"
$ws.Range("E26").Value = "  +6.71%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "15.91"
$ws.Range("E27").Value = "  +1.72%  "
$ws.Range("B28").Value = "Hedera"
$ws.Range("C28").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D28").Value = "0.06845"
$ws.Range("E28").Value = "  +19.12%  "
$ws.Range("D29").Value = "1.346"
$ws.Range("E29").Value = "  +2.96%  "
$ws.Range("E30").Value = "  -0.37%  "
$ws.Range("D31").Value = "3.605"
$ws.Range("E31").Value = "  +2.91%  "
$ws.Range("D32").Value = "3.559"
$ws.Range("E32").Value = "  +2.07%  "
$ws.Range("D33").Value = "1.663"
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("D34").Value = "1.031"
$ws.Range("E34").Value = "  +2.71%  "
$ws.Range("D35").Value = "0.6213"
$ws.Range("E35").Value = "  +4.02%  "
$ws.Range("D36").Value = "2.378"
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("D37").Value = "2.708"
$ws.Range("E37").Value = "  +2.93%  "
$ws.Range("D38").Value = "6.303"
$ws.Range("E38").Value = "  +6.77%  "
$ws.Range("D39").Value = "0.01621"
$ws.Range("E39").Value = "  +1.38%  "
$ws.Range("D40").Value = "1.102.10"
$ws.Range("E40").Value = "  +1.37%  "
$ws.Range("D41").Value = "0.8770"
$ws.Range("E41").Value = "  +1.94%  "
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("E43").Value = "  +1.07%  "
$ws.Range("D44").Value = "1.838.33"
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("D46").Value = "57.12"
$ws.Range("E46").Value = "  +1.98%  "
$ws.Range("E47").Value = "  +0.59%  "
$ws.Range("D48").Value = "8.164"
$ws.Range("E48").Value = "  +1.45%  "
$ws.Range("E49").Value = "  +1.00%  "
$ws.Range("D50").Value = "0.4290"
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("D51").Value = "6.029"
$ws.Range("E51").Value = "  +3.70%  "
